$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B11"  = 6.782000000000001
    "A12"  = -21.526
    "B23"  = 7.628
    "A27"  = -21.637
    "B28"  = 5.196
    "A32"  = -21.23
    "B32"  = 7.417
    "B34"  = 6.739999999999999
    "A36"  = -20.764
    "A38"  = -20.538
    "B42"  = 8.020999999999999
    "A46"  = -21.697
    "B49"  = 6.237
    "A54"  = -21.108
    "B54"  = 5.524
    "A55"  = -22.184
    "A56"  = -21.461
    "A67"  = -21.422
    "A69"  = -21.422
    "A72"  = -21.621
    "B78"  = 7.607000000000001
    "B80"  = 7.665999999999999
    "A83"  = -21.242
    "A86"  = -21.859
    "A91"  = -20.94
    "A93"  = -21.435
    "B97"  = 5.191000000000001
    "A99"  = -21.252
    "B99"  = 5.483000000000001
    "B101" = 5.217000000000001
    "A104" = -21.437
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
